# Apply "update with rough draft" edit to the IDS 590 class schedule.
#
# Summary of content changes (last two class sessions of the schedule):
#   - Row 22 (Groupby class): the In-Class Exercise cell (D22) gains a
#     second bullet linking to the new Opioids Project write-up, in
#     addition to the existing groupby exercise link.
#   - Row 27 (last class, "Thurs, Nov 20"): the old "Machine Learning with
#     scikit-learn" topic/reading/exercise content is replaced with new
#     "Discuss Opioids Project" / "Opioids Working Session" content, and
#     the In-Class Exercise cell (D27) is removed entirely (no exercise
#     for that session), so the row shrinks back to its default height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: append the new Opioids Project link under the existing
#     groupby exercise link, and let the cell wrap like its neighbours.
$ws.Range("D22").Value = '- `Link <../ids720_specific/exercises/exercise_groupby.html>`_' + "`n" + '- `Opioids Project <https://github.com/nickeubank/practicaldatascience_book/blob/main/ids590_specific/opioids_590/PDS590_ProjectSummary.pdf>`_'
$ws.Range("D22").WrapText = $true

# --- Row 27: swap in the Opioids Project session content.
$ws.Range("B27").Value = '- Discuss Opioids Project'
$ws.Range("C27").Value = '- Opioids Working Session'

# Match the formatting used by the other "no in-class exercise this day"
# rows (e.g. D18) by copying its format onto C27, then drop D27 entirely.
$ws.Range("D18").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D27").Clear()

# Row 27 no longer needs the extra height the old multi-line content
# required, so let it size back down to the sheet default.
$ws.Rows.Item(27).AutoFit()

# Restore the previously-selected cell to match the author's final
# cursor position in this revision.
$ws.Range("B27").Select()
